$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.298.64'
$ws.Range('E2').Value = '  +3.75%  '
$ws.Range('D3').Value = '2.627.70'
$ws.Range('E3').Value = '  +3.89%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.11'
$ws.Range('E5').Value = '  +2.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.78'
$ws.Range('E6').Value = '  +1.82%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +1.00%  '
$ws.Range('D9').Value = '2.626.39'
$ws.Range('E9').Value = '  +3.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.166'
$ws.Range('E10').Value = '  +13.63%  '
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('E12').Value = '  +2.22%  '
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.114.19'
$ws.Range('E14').Value = '  +4.12%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000189'
$ws.Range('E15').Value = '  +10.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.70'
$ws.Range('E16').Value = '  +1.71%  '
$ws.Range('D17').Value = '71.207.27'
$ws.Range('E17').Value = '  +4.26%  '
$ws.Range('D18').Value = '2.629.73'
$ws.Range('E18').Value = '  +4.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '381.51'
$ws.Range('E19').Value = '  +8.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.91'
$ws.Range('E20').Value = '  +5.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.52'
$ws.Range('E21').Value = '  +3.53%  '
$ws.Range('E22').Value = '  -1.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.56'
$ws.Range('E23').Value = '  +2.21%  '
$ws.Range('E24').Value = '  +5.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.97'
$ws.Range('E25').Value = '  +15.49%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.71'
$ws.Range('E27').Value = '  +7.67%  '
$ws.Range('D28').Value = '2.764.53'
$ws.Range('E28').Value = '  +2.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').Value = '0.0₃0967'
$ws.Range('E30').Value = '  +7.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '548.31'
$ws.Range('E31').Value = '  +7.57%  '
$ws.Range('E32').Value = '  +3.19%  '
$ws.Range('E33').Value = '  +5.45%  '
$ws.Range('E34').Value = '  +3.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '166.08'
$ws.Range('E36').Value = '  +1.90%  '
$ws.Range('E37').Value = '  -1.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.24'
$ws.Range('E38').Value = '  +4.42%  '
$ws.Range('E39').Value = '  +6.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.03'
$ws.Range('E40').Value = '  +1.87%  '
$ws.Range('E41').Value = '  +4.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.64'
$ws.Range('E42').Value = '  +9.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.04'
$ws.Range('E44').Value = '  +4.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.331'
$ws.Range('E45').Value = '  +1.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.97'
$ws.Range('E46').Value = '  +2.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '154.53'
$ws.Range('E47').Value = '  +0.69%  '
$ws.Range('E48').Value = '  +1.90%  '
$ws.Range('E49').Value = '  +5.26%  '
$ws.Range('E50').Value = '  +2.49%  '
$ws.Range('D51').Value = '0.0₆0264'
$ws.Range('E51').Value = '  +1.62%  '
